# Adds a new "FirstParagraph"-styled paragraph right after the author
# line ("John T. Foster") and right before the bookmarks that wrap the
# students-graduated table. The paragraph talks about the author's blog
# and contains a hyperlink ("scientific computing topics") styled as
# InternetLink.

$d = $word.ActiveDocument

# --- 1. Create a brand new (empty) paragraph right after the author
#        line. Using InsertParagraphAfter (rather than a Find/Replace
#        that rewrites "Foster" itself) leaves the existing run
#        structure of the author paragraph completely untouched. -----
$authorPara = $d.Paragraphs(2).Range
$authorPara.InsertParagraphAfter()

$newPara = $d.Paragraphs(3)
$newPara.Range.Style = "FirstParagraph"

# --- 2. Fill in the paragraph's text ----------------------------------
$blogText = "With regard to assisting (primarily graduate) students in their learning process, I maintain an active blog of helpful tips and answered questions. Generally these are related to scientific computing topics. I try to maintain a philosophy where anytime I believe a student’s question will come up again in the future, I write a blog detailing the path forward as opposed to providing an individual answer. This has paid off numerous times where I can simply point the students to a blog post where I’ve already answered their specific or a similar question."
$newPara.Range.Text = $blogText

# --- 3. Turn "scientific computing topics" into a hyperlink ----------
$hlRange = $newPara.Range.Duplicate
$hlRange.Find.Execute("scientific computing topics", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$d.Hyperlinks.Add($hlRange, "https://johntfoster.github.io/blog/")

# Restyle the newly created hyperlink run so it matches the document's
# InternetLink character style (Hyperlinks.Add defaults to "Hyperlink").
$hlRange2 = $newPara.Range.Duplicate
$hlRange2.Find.Execute("scientific computing topics", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$hlRange2.Style = "InternetLink"

Write-Output "Inserted mentoring-blog paragraph with hyperlink."
